$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 77.375
$ws.Range("I8").Value = 77.375
$ws.Range("K8").Value = 232.125
$ws.Range("M8").Value = -93.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 467.17648
$ws.Range("I28").Value = 478
$ws.Range("J28").Value = 416.66666
$ws.Range("K28").Value = 478
$ws.Range("L28").Value = 416.66666
$ws.Range("M28").Value = 7
$ws.Range("N28").Value = -1386.66666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 52782
$ws.Range("I33").Value = 28731.828
$ws.Range("J33").Value = 333367.34
$ws.Range("K33").Value = 28731.828
$ws.Range("L33").Value = 333367.34
$ws.Range("M33").Value = -28502.828
$ws.Range("N33").Value = -333825.34

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2346.1936
$ws.Range("I138").Value = 1131.1154
$ws.Range("J138").Value = 3223.75
$ws.Range("K138").Value = 3393.3462
$ws.Range("L138").Value = 9671.25
$ws.Range("M138").Value = 1746.6538
$ws.Range("N138").Value = -19951.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1288.4615
$ws.Range("I2").Value = 1288.4615
$ws.Range("K2").Value = 1288.4615
$ws.Range("M2").Value = -1175.4615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1673.8
$ws.Range("I45").Value = 1790.5
$ws.Range("J45").Value = 1207
$ws.Range("K45").Value = 1790.5
$ws.Range("L45").Value = 1207
$ws.Range("M45").Value = -1413.5
$ws.Range("N45").Value = -1961

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 26330
$ws.Range("J86").Value = 26330
$ws.Range("L86").Value = 26330
$ws.Range("N86").Value = -28702

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 26330
$ws.Range("J89").Value = 26330
$ws.Range("L89").Value = 78990
$ws.Range("N89").Value = -90846

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 23120
$ws.Range("J109").Value = 23120
$ws.Range("L109").Value = 23120
$ws.Range("N109").Value = -25894

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1288.4615
$ws.Range("I116").Value = 1288.4615
$ws.Range("K116").Value = 1288.4615
$ws.Range("M116").Value = 1005.5385

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 93212
$ws.Range("J139").Value = 93212
$ws.Range("L139").Value = 93212
$ws.Range("N139").Value = -103492

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1288.4615
$ws.Range("I3").Value = 1288.4615
$ws.Range("K3").Value = 1288.4615
$ws.Range("M3").Value = -1174.4615

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2852.7144
$ws.Range("I86").Value = 2852.7144
$ws.Range("K86").Value = 2852.7144
$ws.Range("M86").Value = -1729.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2852.7144
$ws.Range("I89").Value = 2852.7144
$ws.Range("K89").Value = 14263.572
$ws.Range("M89").Value = -8647.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 716.7241
$ws.Range("I94").Value = 711.03705
$ws.Range("J94").Value = 793.5
$ws.Range("K94").Value = 711.03705
$ws.Range("L94").Value = 793.5
$ws.Range("M94").Value = -260.03705
$ws.Range("N94").Value = -1695.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1022.6
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 1506.5
$ws.Range("K107").Value = 700
$ws.Range("L107").Value = 1506.5
$ws.Range("M107").Value = 1220
$ws.Range("N107").Value = -5346.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 705.44446
$ws.Range("I16").Value = 674.8333
$ws.Range("J16").Value = 766.6667
$ws.Range("K16").Value = 674.8333
$ws.Range("L16").Value = 766.6667
$ws.Range("M16").Value = -387.8333
$ws.Range("N16").Value = -1340.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2000
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 2000
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 43500
$ws.Range("I52").Value = 45000
$ws.Range("J52").Value = 42750
$ws.Range("K52").Value = 45000
$ws.Range("L52").Value = 42750
$ws.Range("M52").Value = -44706
$ws.Range("N52").Value = -43338

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 705.44446
$ws.Range("I113").Value = 674.8333
$ws.Range("J113").Value = 766.6667
$ws.Range("K113").Value = 674.8333
$ws.Range("L113").Value = 766.6667
$ws.Range("M113").Value = 1495.1667
$ws.Range("N113").Value = -5106.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 327.83334
$ws.Range("I2").Value = 17.25
$ws.Range("J2").Value = 416.57144
$ws.Range("K2").Value = 103.5
$ws.Range("L2").Value = 2499.42864
$ws.Range("M2").Value = 9.5
$ws.Range("N2").Value = -2725.42864

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1729.3334
$ws.Range("I113").Value = 2553.1428
$ws.Range("J113").Value = 576
$ws.Range("K113").Value = 7659.428400000001
$ws.Range("L113").Value = 1728
$ws.Range("M113").Value = -5489.428400000001
$ws.Range("N113").Value = -6068

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 934.35297
$ws.Range("I122").Value = 403.75
$ws.Range("J122").Value = 1406
$ws.Range("K122").Value = 3633.75
$ws.Range("L122").Value = 12654
$ws.Range("M122").Value = -1183.75
$ws.Range("N122").Value = -17554

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 7249.778
$ws.Range("J9").Value = 15727
$ws.Range("L9").Value = 15727
$ws.Range("N9").Value = -16067

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5117.3335
$ws.Range("I70").Value = 4945.9443
$ws.Range("K70").Value = 4945.9443
$ws.Range("M70").Value = -4675.9443

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5117.3335
$ws.Range("I73").Value = 4945.9443
$ws.Range("K73").Value = 4945.9443
$ws.Range("M73").Value = -4009.9443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 766
$ws.Range("I93").Value = 607.1429
$ws.Range("J93").Value = 951.3333
$ws.Range("K93").Value = 607.1429
$ws.Range("L93").Value = 951.3333
$ws.Range("M93").Value = 640.8571
$ws.Range("N93").Value = -3447.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2614.8684
$ws.Range("I100").Value = 1424.375
$ws.Range("J100").Value = 2932.3333
$ws.Range("K100").Value = 1424.375
$ws.Range("L100").Value = 2932.3333
$ws.Range("M100").Value = -883.375
$ws.Range("N100").Value = -4014.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 31666.666
$ws.Range("J19").Value = 31666.666
$ws.Range("L19").Value = 31666.666
$ws.Range("N19").Value = -32014.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 10919.625
$ws.Range("I43").Value = 3887.8333
$ws.Range("K43").Value = 3887.8333
$ws.Range("M43").Value = -3738.8333

Write-Host "Done applying market price updates"